# Update TTLE logit exponents for nonroad (aircraft, rail, ships) from -1 to -0.1
# to avoid floating point overflow with very high tech costs, and update the
# "About" sheet notes to explain the calibration / non-road exception.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. TTLE sheet: change nonroad (aircraft/rail/ships) logit exponents
# ---------------------------------------------------------------------
$ttle = $wb.Worksheets.Item("TTLE")

$ttle.Range("B4:C6").Value = -0.1
$ttle.Range("B4:C6").Select() | Out-Null

# ---------------------------------------------------------------------
# 2. About sheet: rework the Notes section
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Replace row 12's note about PNNL's GCAM model with two new notes about
# calibration of onroad sectors vs. the nonroad exception.
# (Written in this order so new shared strings land in the same order as
# the target workbook's sharedStrings.xml table.)
$about.Range("A13").Value = "For non-road we use -0.1 because of costs preventing the model from solving."
$about.Range("A12").Value = "We use calibrated values in onroad sectors."

# Row 15/16 used to read "For more on this, see the ""Modified Logit""..."
# followed by the URL. The wording changes to "Unmodified Logit" and the
# URL moves down one row to make room for the new note above.
$about.Range("A15").Value = 'For more on this, see the "Unmodified Logit" equation description at:'
$about.Range("A16").Value = "https://jgcri.github.io/gcam-doc/choice.html"

# ---------------------------------------------------------------------
# 3. Selection / active-tab bookkeeping to mirror the saved workbook state
# ---------------------------------------------------------------------
$about.Activate() | Out-Null
$about.Range("A16").Select() | Out-Null
